# "Implemented Label wise bill generate(Individual)"
#
# Fills in the teacher's identity on the bill header (name / position /
# department), records the individual quantities for three bill lines
# (question-paper-setting, answer-script-examining and invigilation-style
# rows), and states the grand total in words. The dependent formula cells
# (I16, I20, I26, I32) recompute automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header block: teacher's name / designation / department.
$ws.Range("A3").Value  = "নাম: Dr. Muhammad Sheikh Sadi"
$ws.Range("A4").Value  = "পদবী: অধ্যাপক"
$ws.Range("F5").Value  = "বিভাগ :সিএসই"

# Bill-line quantities.
$ws.Range("G16").Value = 27
$ws.Range("G20").Value = 40
$ws.Range("G26").Value = 1

# Grand total, spelled out in words.
$ws.Range("A32").Value = "কথায়:সাত হাজার বাহান্ন টাকা মাত্র।"

# Leave the selection where the author left it.
$ws.Range("B5").Select() | Out-Null
